$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.416.04"
$ws.Range("E2").Value = "  -1.40%  "
$ws.Range("D3").Value = "2.636.95"
$ws.Range("E3").Value = "  -2.75%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.29"
$ws.Range("E5").Value = "  -2.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.71"
$ws.Range("E6").Value = "  +1.18%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -3.54%  "
$ws.Range("D9").Value = "2.635.15"
$ws.Range("E9").Value = "  -2.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.143"
$ws.Range("E10").Value = "  -1.46%  "
$ws.Range("E11").Value = "  +1.27%  "
$ws.Range("E12").Value = "  -0.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.26"
$ws.Range("E13").Value = "  -0.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.89"
$ws.Range("E14").Value = "  -1.82%  "
$ws.Range("E15").Value = "  -2.90%  "
$ws.Range("E16").Value = "  -2.64%  "
$ws.Range("D17").Value = "67.306.88"
$ws.Range("E17").Value = "  -1.43%  "
$ws.Range("D18").Value = "2.598.29"
$ws.Range("E18").Value = "  -4.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.19"
$ws.Range("E19").Value = "  +3.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "360.82"
$ws.Range("E21").Value = "  -2.42%  "
$ws.Range("E22").Value = "  -2.58%  "
$ws.Range("E23").Value = "  -4.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.93"
$ws.Range("E24").Value = "  +9.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.96"
$ws.Range("E25").Value = "  -5.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "70.50"
$ws.Range("E27").Value = "  -3.46%  "
$ws.Range("D28").Value = "2.769.36"
$ws.Range("E28").Value = "  -2.71%  "
$ws.Range("E29").Value = "  -2.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "556.86"
$ws.Range("E31").Value = "  -3.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.94"
$ws.Range("E32").Value = "  -2.14%  "
$ws.Range("E33").Value = "  -2.92%  "
$ws.Range("E34").Value = "  -3.44%  "
$ws.Range("E35").Value = "  +3.78%  "
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("E37").Value = "  -4.81%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "157.38"
$ws.Range("E38").Value = "  -1.88%  "
$ws.Range("E39").Value = "  -3.28%  "
$ws.Range("E40").Value = "  -2.75%  "
$ws.Range("E41").Value = "  -2.89%  "
$ws.Range("E42").Value = "  -3.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.95"
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("E45").Value = "  -4.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.22"
$ws.Range("E46").Value = "  -1.29%  "
$ws.Range("E47").Value = "  -3.17%  "
$ws.Range("E48").Value = "  -1.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "152.43"
$ws.Range("E49").Value = "  -1.67%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.82"
$ws.Range("E50").Value = "  -1.85%  "
$ws.Range("E51").Value = "  -1.35%  "
